$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a string to a cell while forcing text storage (so values
# that look numeric, e.g. "580.16" or "25.00", stay text exactly like the
# original inlineStr cells and aren't coerced into numbers). A leading
# apostrophe is Excel's standard "treat as text" prefix and is stripped from
# the stored value automatically.
function Set-TextValue($addr, $value) {
    $ws.Range($addr).Value = "'" + $value
}

# Row 2 (Bitcoin)
Set-TextValue "D2" "67.091.93"
$ws.Range("E2").Value = "  +0.32%  "

# Row 3 (Ethereum)
Set-TextValue "D3" "3.122.75"
$ws.Range("E3").Value = "  +0.64%  "

# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  -0.03%  "

# Row 5 (BNB)
Set-TextValue "D5" "580.16"
$ws.Range("E5").Value = "  -0.22%  "

# Row 6 (Solana)
Set-TextValue "D6" "173.98"
$ws.Range("E6").Value = "  +0.44%  "

# Row 7 (USDC)
Set-TextValue "D7" "0.999"
$ws.Range("E7").Value = "  -0.08%  "

# Row 8 (XRP)
$ws.Range("E8").Value = "  -0.34%  "

# Row 9 (now Dogecoin, was Toncoin)
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue "D9" "0.155"
$ws.Range("E9").Value = "  -0.62%  "

# Row 10 (now Toncoin, was Dogecoin)
$ws.Range("B10").Value = "Toncoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D10" "6.41"
$ws.Range("E10").Value = "  -0.50%  "

# Row 11 (Cardano)
Set-TextValue "D11" "0.481"
$ws.Range("E11").Value = "  -0.79%  "

# Row 12 (ShibaInu)
$ws.Range("E12").Value = "  -0.11%  "

# Row 13 (Avalanche)
Set-TextValue "D13" "37.25"
$ws.Range("E13").Value = "  -0.47%  "

# Row 14 (TRON)
$ws.Range("E14").Value = "  -1.68%  "

# Row 15 (WrappedliquidstakedEther2.0)
Set-TextValue "D15" "3.639.85"
$ws.Range("E15").Value = "  +0.69%  "

# Row 16 (WrappedBTC)
Set-TextValue "D16" "67.069.70"
$ws.Range("E16").Value = "  +0.24%  "

# Row 17 (Polkadot)
$ws.Range("E17").Value = "  -0.76%  "

# Row 18 (WrappedEther)
Set-TextValue "D18" "3.122.54"
$ws.Range("E18").Value = "  +0.47%  "

# Row 19 (Chainlink)
Set-TextValue "D19" "16.43"
$ws.Range("E19").Value = "  +2.03%  "

# Row 20 (BitcoinCash)
Set-TextValue "D20" "491.72"
$ws.Range("E20").Value = "  +2.09%  "

# Row 21 (Uniswap)
$ws.Range("E21").Value = "  +5.82%  "

# Row 22 (Polygon)
Set-TextValue "D22" "0.708"
$ws.Range("E22").Value = "  -0.98%  "

# Row 23 (Litecoin)
Set-TextValue "D23" "84.17"
$ws.Range("E23").Value = "  +0.15%  "

# Row 24 (InternetComputer(DFINITY))
$ws.Range("E24").Value = "  +0.70%  "

# Row 25 (Fetch.AI)
$ws.Range("E25").Value = "  -3.42%  "

# Row 27 (Dai)
$ws.Range("E27").Value = "  -0.05%  "

# Row 28 (NEARProtocol)
$ws.Range("E28").Value = "  -1.06%  "

# Row 29 (ImmutableX)
$ws.Range("E29").Value = "  -1.67%  "

# Row 30 (PancakeSwap)
$ws.Range("E30").Value = "  -0.54%  "

# Row 31 (EthereumClassic)
Set-TextValue "D31" "28.65"
$ws.Range("E31").Value = "  -0.19%  "

# Row 32 (Hedera)
$ws.Range("E32").Value = "  -0.48%  "

# Row 33 (PEPE)
Set-TextValue "D33" ([string]::Concat("0.0", [char]8323, "0949"))
$ws.Range("E33").Value = "  -6.09%  "

# Row 34 (FirstDigitalUSD)
$ws.Range("E34").Value = "  -0.11%  "

# Row 35 (Filecoin)
Set-TextValue "D35" "5.89"
$ws.Range("E35").Value = "  -0.43%  "

# Row 37 (Arweave)
Set-TextValue "D37" "47.34"
$ws.Range("E37").Value = "  -1.72%  "

# Row 38 (Stacks)
Set-TextValue "D38" "2.05"
$ws.Range("E38").Value = "  -3.05%  "

# Row 39 (TheGraph)
$ws.Range("E39").Value = "  -2.36%  "

# Row 40 (Kaspa)
$ws.Range("E40").Value = "  +1.50%  "

# Row 41 (Cosmos)
Set-TextValue "D41" "8.54"
$ws.Range("E41").Value = "  -1.50%  "

# Row 42 (Maker)
Set-TextValue "D42" "2.824.86"
$ws.Range("E42").Value = "  -0.13%  "

# Row 43 (now Bittensor, was dogwifhat)
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D43" "383.98"
$ws.Range("E43").Value = "  -0.01%  "

# Row 44 (now dogwifhat, was Bittensor)
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D44" "2.61"
$ws.Range("E44").Value = "  -7.54%  "

# Row 45 (VeChain)
$ws.Range("E45").Value = "  -2.38%  "

# Row 46 (Monero)
$ws.Range("E46").Value = "  +0.64%  "

# Row 47 (USDe)
$ws.Range("E47").Value = "  -0.01%  "

# Row 48 (InjectiveProtocol)
Set-TextValue "D48" "25.00"
$ws.Range("E48").Value = "  +0.43%  "

# Row 49 (ThetaToken)
$ws.Range("E49").Value = "  -0.97%  "

# Row 50 (Stellar)
$ws.Range("E50").Value = "  -0.73%  "

# Row 51 (THORChain)
Set-TextValue "D51" "6.76"
$ws.Range("E51").Value = "  -0.94%  "
